$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the Seattle row for the "Top 100 zip" large-dataset
# statistics, and carry over the same formatting used by the row above
# (the "JSON flat files" header row) so borders/fonts match the table.
$ws.Rows(5).Insert()
$ws.Range("A4:I4").Copy()
$ws.Range("A5:I5").PasteSpecial(-4122)

$ws.Range("A5").Value = "Top 100 zip"
$ws.Range("B5").Value = 1669
$ws.Range("C5").Value = 3973
$ws.Range("D5").Value = 2567
$ws.Range("E5").Value = 1711
$ws.Range("F5").Value = 60
$ws.Range("G5").Value = 40
$ws.Range("H5").Value = 0.044914580000000003
$ws.Range("I5").Value = 0.033973719999999999

# Footnote explaining the gender association caveat for the new numbers.
$ws.Range("C9").Value = "*with male/female gender association"
$ws.Range("C9").Font.Size = 8

# Relabel the reviews column header to point at the footnote.
$ws.Range("C3").Value = "# Reviews*"

# Column A needed to widen slightly to fit "Top 100 zip".
$ws.Columns("A").ColumnWidth = 14.77734375

# Leave the selection where the author left it after finishing the edits.
$ws.Range("B14").Select()
